$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3547.1191
$ws.Range("I40").Value = 3017.9375
$ws.Range("K40").Value = 3017.9375
$ws.Range("M40").Value = -2842.9375
$ws.Range("H41").Value = 4448.6665
$ws.Range("J41").Value = 3500.5
$ws.Range("L41").Value = 3500.5
$ws.Range("N41").Value = -4380.5
$ws.Range("H43").Value = 5057.8
$ws.Range("I43").Value = 4225.3335
$ws.Range("J43").Value = 5320.684
$ws.Range("K43").Value = 4225.3335
$ws.Range("L43").Value = 5320.684
$ws.Range("M43").Value = -4156.3335
$ws.Range("N43").Value = -5458.684
$ws.Range("H86").Value = 4441
$ws.Range("I86").Value = 3890.5715
$ws.Range("J86").Value = 5083.1665
$ws.Range("K86").Value = 3890.5715
$ws.Range("L86").Value = 5083.1665
$ws.Range("M86").Value = -2767.5715
$ws.Range("N86").Value = -7329.1665
$ws.Range("H89").Value = 4441
$ws.Range("I89").Value = 3890.5715
$ws.Range("J89").Value = 5083.1665
$ws.Range("K89").Value = 19452.8575
$ws.Range("L89").Value = 25415.8325
$ws.Range("M89").Value = -13836.8575
$ws.Range("N89").Value = -36647.8325
$ws.Range("H106").Value = 7401.206
$ws.Range("I106").Value = 2755.8635
$ws.Range("K106").Value = 2755.8635
$ws.Range("M106").Value = -2124.8635
$ws.Range("H129").Value = 1690.7778
$ws.Range("J129").Value = 3089.6667
$ws.Range("L129").Value = 9269.000100000001
$ws.Range("N129").Value = -19269.0001
$ws.Range("H131").Value = 252961.4
$ws.Range("I131").Value = 2386
$ws.Range("J131").Value = 1255263
$ws.Range("K131").Value = 7158
$ws.Range("L131").Value = 3765789
$ws.Range("M131").Value = -2118
$ws.Range("N131").Value = -3775869
$ws.Range("H138").Value = 4277.1875
$ws.Range("I138").Value = 3649.077
$ws.Range("J138").Value = 6999
$ws.Range("K138").Value = 10947.231
$ws.Range("L138").Value = 20997
$ws.Range("M138").Value = -5807.231
$ws.Range("N138").Value = -31277

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4175.8965
$ws.Range("I132").Value = 3380.9546
$ws.Range("K132").Value = 10142.8638
$ws.Range("M132").Value = -7612.863799999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1659
$ws.Range("I99").Value = 1629
$ws.Range("K99").Value = 1629
$ws.Range("M99").Value = -131
$ws.Range("H105").Value = 13789.708
$ws.Range("I105").Value = 13827.412
$ws.Range("J105").Value = 13698.143
$ws.Range("K105").Value = 13827.412
$ws.Range("L105").Value = 13698.143
$ws.Range("M105").Value = -12080.412
$ws.Range("N105").Value = -17192.143
$ws.Range("H134").Value = 2455.5625
$ws.Range("I134").Value = 1598.3077
$ws.Range("K134").Value = 4794.9231
$ws.Range("M134").Value = -2259.9231

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38417.066
$ws.Range("I31").Value = 4553.4736
$ws.Range("K31").Value = 4553.4736
$ws.Range("M31").Value = -4258.4736
$ws.Range("H34").Value = 38417.066
$ws.Range("I34").Value = 4553.4736
$ws.Range("K34").Value = 4553.4736
$ws.Range("M34").Value = -4351.4736
$ws.Range("H58").Value = 3432.742
$ws.Range("I58").Value = 1803.68
$ws.Range("K58").Value = 1803.68
$ws.Range("M58").Value = -1600.68
$ws.Range("H86").Value = 14803
$ws.Range("I86").Value = 11251.75
$ws.Range("K86").Value = 11251.75
$ws.Range("M86").Value = -10128.75
$ws.Range("H89").Value = 14803
$ws.Range("I89").Value = 11251.75
$ws.Range("K89").Value = 56258.75
$ws.Range("M89").Value = -50642.75
$ws.Range("H99").Value = 1989.6072
$ws.Range("I99").Value = 1801.762
$ws.Range("J99").Value = 2553.1428
$ws.Range("K99").Value = 1801.762
$ws.Range("L99").Value = 2553.1428
$ws.Range("M99").Value = -303.7619999999999
$ws.Range("N99").Value = -5549.1428
$ws.Range("H126").Value = 1989.6072
$ws.Range("I126").Value = 1801.762
$ws.Range("J126").Value = 2553.1428
$ws.Range("K126").Value = 5405.286
$ws.Range("L126").Value = 7659.428400000001
$ws.Range("M126").Value = -2935.286
$ws.Range("N126").Value = -12599.4284
$ws.Range("H132").Value = 2134.524
$ws.Range("I132").Value = 1267.2778
$ws.Range("K132").Value = 3801.8334
$ws.Range("M132").Value = -1271.8334
$ws.Range("H134").Value = 2686.125
$ws.Range("I134").Value = 1672.7646
$ws.Range("K134").Value = 5018.293799999999
$ws.Range("M134").Value = -2483.293799999999
$ws.Range("H136").Value = 3432.742
$ws.Range("I136").Value = 1803.68
$ws.Range("K136").Value = 5411.04
$ws.Range("M136").Value = -2861.04

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2822.2
$ws.Range("J34").Value = 3925
$ws.Range("L34").Value = 11775
$ws.Range("N34").Value = -11943
$ws.Range("H39").Value = 2504.4443
$ws.Range("I39").Value = 2128
$ws.Range("J39").Value = 2975
$ws.Range("K39").Value = 6384
$ws.Range("L39").Value = 8925
$ws.Range("M39").Value = -6090
$ws.Range("N39").Value = -9513
$ws.Range("H55").Value = 1662.4
$ws.Range("I55").Value = 1404
$ws.Range("J55").Value = 2050
$ws.Range("K55").Value = 4212
$ws.Range("L55").Value = 6150
$ws.Range("M55").Value = -4035
$ws.Range("N55").Value = -6504
$ws.Range("H68").Value = 2574.75
$ws.Range("J68").Value = 2649.5
$ws.Range("L68").Value = 7948.5
$ws.Range("N68").Value = -9570.5
$ws.Range("H71").Value = 2574.75
$ws.Range("J71").Value = 2649.5
$ws.Range("L71").Value = 23845.5
$ws.Range("N71").Value = -31957.5
$ws.Range("H81").Value = 3995
$ws.Range("I81").Value = 1566.6666
$ws.Range("K81").Value = 4699.9998
$ws.Range("M81").Value = -3576.9998
$ws.Range("H84").Value = 3995
$ws.Range("I84").Value = 1566.6666
$ws.Range("K84").Value = 14099.9994
$ws.Range("M84").Value = -8483.999400000001
$ws.Range("H112").Value = 100010810
$ws.Range("I112").Value = 125006010
$ws.Range("J112").Value = 30000
$ws.Range("K112").Value = 375018030
$ws.Range("L112").Value = 90000
$ws.Range("M112").Value = -375016922
$ws.Range("N112").Value = -92216
$ws.Range("H122").Value = 1676.0435
$ws.Range("J122").Value = 1676.0435
$ws.Range("L122").Value = 15084.3915
$ws.Range("N122").Value = -19984.3915
$ws.Range("H133").Value = 7163.385
$ws.Range("I133").Value = 3418
$ws.Range("K133").Value = 10254
$ws.Range("M133").Value = -5194

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 11030
$ws.Range("I52").Value = 11030
$ws.Range("K52").Value = 11030
$ws.Range("M52").Value = -10771
$ws.Range("H70").Value = 11532.083
$ws.Range("I70").Value = 6637.6
$ws.Range("K70").Value = 6637.6
$ws.Range("M70").Value = -6367.6
$ws.Range("H73").Value = 11532.083
$ws.Range("I73").Value = 6637.6
$ws.Range("K73").Value = 6637.6
$ws.Range("M73").Value = -5701.6
$ws.Range("H132").Value = 44627.4
$ws.Range("I132").Value = 49944.137
$ws.Range("J132").Value = 5638
$ws.Range("K132").Value = 149832.411
$ws.Range("L132").Value = 16914
$ws.Range("M132").Value = -147302.411
$ws.Range("N132").Value = -21974

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 40000
$ws.Range("I42").Value = 40000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 40000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -39437
$ws.Range("N42").ClearContents()
$ws.Range("H49").Value = 40000
$ws.Range("I49").Value = 40000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 40000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -39853
$ws.Range("N49").ClearContents()
$ws.Range("H55").Value = 1924196.5
$ws.Range("I55").Value = 3125266.8
$ws.Range("J55").Value = 2484.1
$ws.Range("K55").Value = 3125266.8
$ws.Range("L55").Value = 2484.1
$ws.Range("M55").Value = -3125093.8
$ws.Range("N55").Value = -2830.1
$ws.Range("H132").Value = 3598.3333
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 9293.333000000001
$ws.Range("K132").Value = 5100
$ws.Range("L132").Value = 27879.999
$ws.Range("M132").Value = -2570
$ws.Range("N132").Value = -32939.999
$ws.Range("H136").Value = 9119.177
$ws.Range("I136").Value = 2446.4443
$ws.Range("K136").Value = 7339.3329
$ws.Range("M136").Value = -4789.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 11500
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 11500
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H102").Value = 74992.5
$ws.Range("J102").Value = 74992.5
$ws.Range("L102").Value = 74992.5
$ws.Range("N102").Value = -81482.5
$ws.Range("H122").Value = 7238.3237
$ws.Range("I122").Value = 3072.5
$ws.Range("J122").Value = 11924.875
$ws.Range("K122").Value = 9217.5
$ws.Range("L122").Value = 35774.625
$ws.Range("M122").Value = -6767.5
$ws.Range("N122").Value = -40674.625
$ws.Range("H123").Value = 60429
$ws.Range("J123").Value = 60429
$ws.Range("L123").Value = 60429
$ws.Range("N123").Value = -70229
$ws.Range("H132").Value = 6124.9243
$ws.Range("I132").Value = 5481.898
$ws.Range("K132").Value = 16445.694
$ws.Range("M132").Value = -13915.694
$ws.Range("H136").Value = 6874.25
$ws.Range("I136").Value = 3998.4443
$ws.Range("K136").Value = 11995.3329
$ws.Range("M136").Value = -9445.332900000001
